$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Falsche Zeitangabe soll eine erneute Eingabe erzwingen: Status auf "ok" setzen
$ws.Range("C16").Value = "ok"
$ws.Range("C17").Value = "ok"

# Ansicht: auf den neuen Status-Bereich scrollen und C18 auswählen
$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 7
